$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.975.33'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.76%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.621.52'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.35'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '166.16'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.532'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.621.78'
$ws.Range('D9').Style = "Normal"
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.358'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.22'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.63'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.115.14'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.39%  '
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.010.68'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.622.20'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.11'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.15%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.06'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.93%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '356.61'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.29%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.31'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.27%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.65'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +9.87%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.92'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -7.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '70.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.757.64'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '550.71'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('E34').Value = '  -3.40%  '
$ws.Range('E35').Value = '  +4.05%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  -5.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '155.31'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.07'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.70%  '
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('E41').Value = '  -3.55%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.14'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.42%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.93'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.42'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.24'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('E47').Value = '  -2.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.579'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '151.25'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.78'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.71'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.31%  '
